$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values entered in the author's original order (this determines the
# order new entries land in xl/sharedStrings.xml), even though the rows end up
# sorted by ppm_i in the final worksheet.
$ws.Range("A6").Value = "Lactate"
$ws.Range("A10").Value = "Alanine"
$ws.Range("A24").Value = "Unknow"
$ws.Range("A7").Value = "Unknow_1"
$ws.Range("A8").Value = "Unknow_2"
$ws.Range("A9").Value = "Unknow_3"
$ws.Range("A11").Value = "2-Aminoisobutyric acid"
$ws.Range("A12").Value = "Unknow_4"
$ws.Range("A13").Value = "L-arginine/Leucine"
$ws.Range("A14").Value = "Unknow_5"
$ws.Range("A15").Value = "Unknow_6"
$ws.Range("A16").Value = "Acetate"
$ws.Range("A17").Value = "Acetylphosphate"
$ws.Range("A18").Value = "Homoserine"
$ws.Range("A19").Value = "Acetylglycine"
$ws.Range("A20").Value = "Acetylcysteine"
$ws.Range("A21").Value = "Acetylcholine"
$ws.Range("A22").Value = "Acetone/3-hidroxybutyrate"
$ws.Range("A23").Value = "Piruvate"
$ws.Range("A25").Value = "Beta-alanine"

# ppm_i / ppm_f values for the new rows (6-25), row by row.
$ws.Range("B6").Value = 1.32
$ws.Range("C6").Value = 1.36

$ws.Range("B7").Value = 1.3680000000000001
$ws.Range("C7").Value = 1.38

$ws.Range("B8").Value = 1.38
$ws.Range("C8").Value = 1.405

$ws.Range("B9").Value = 1.425
$ws.Range("C9").Value = 1.4350000000000001

$ws.Range("B10").Value = 1.4650000000000001
$ws.Range("C10").Value = 1.5049999999999999

$ws.Range("B11").Value = 1.5
$ws.Range("C11").Value = 1.5149999999999999

$ws.Range("B12").Value = 1.5149999999999999
$ws.Range("C12").Value = 1.53

$ws.Range("B13").Value = 1.66
$ws.Range("C13").Value = 1.8

$ws.Range("B14").Value = 1.81
$ws.Range("C14").Value = 1.86

$ws.Range("B15").Value = 1.865
$ws.Range("C15").Value = 1.88

$ws.Range("B16").Value = 1.893
$ws.Range("C16").Value = 1.915

$ws.Range("B17").Value = 1.915
$ws.Range("C17").Value = 1.94

$ws.Range("B18").Value = 1.98
$ws.Range("C18").Value = 2.0299999999999998

$ws.Range("B19").Value = 2.0499999999999998
$ws.Range("C19").Value = 2.0649999999999999

$ws.Range("B20").Value = 2.0659999999999998
$ws.Range("C20").Value = 2.077

$ws.Range("B21").Value = 2.13
$ws.Range("C21").Value = 2.15

$ws.Range("B22").Value = 2.2349999999999999
$ws.Range("C22").Value = 2.3050000000000002

$ws.Range("B23").Value = 2.4
$ws.Range("C23").Value = 2.42

$ws.Range("B24").Value = 2.4900000000000002
$ws.Range("C24").Value = 2.54

$ws.Range("B25").Value = 2.5499999999999998
$ws.Range("C25").Value = 2.585

# Column A width (closest achievable quantum to the target 18.5546875 OOXML width).
$ws.Columns.Item(1).ColumnWidth = 17.71

# Final selection, matching the author's last-saved cursor position.
$null = $ws.Range("A25").Select()
